# LTIM.NS.xlsx weekly-data update — "break out stock.yaml completed"
#
# 1) Q56 corrected from 2 -> 0
# 2) O413 corrected from 0 -> 2
# 3) R415 / R416 (previously blank) filled in with 0
# 4) 8 new weekly rows (417-424, 2024-07-01 .. 2024-08-19) appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- point fixes on existing rows -------------------------------------
$ws.Cells.Item(56, 17).Value = 0     # Q56: 2 -> 0
$ws.Cells.Item(413, 15).Value = 2    # O413: 0 -> 2
$ws.Cells.Item(415, 18).Value = 0    # R415: blank -> 0
$ws.Cells.Item(416, 18).Value = 0    # R416: blank -> 0

# --- new rows ----------------------------------------------------------
# Columns: A Datetime, B Open, C High, D Low, E Close, F Adj Close,
#          G Volume, H Year, I Month, J Day, K Hour, L Minute, M Second,
#          N Week, O isPivot, P two_line_structure, Q detect_structure
$newRows = @(
    @(417, 45474, 5372.4501953125, 5550, 5352.14990234375, 5421.7001953125, 5421.7001953125, 2314757, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
    @(418, 45481, 5421.7001953125, 5659.7001953125, 5320.35009765625, 5572.64990234375, 5572.64990234375, 2181606, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(419, 45488, 5612, 5879.14990234375, 5438.2998046875, 5762.75, 5762.75, 4215383, 2024, 7, 15, 0, 0, 0, 29, 1, 0, 1),
    @(420, 45495, 5733, 5811.89990234375, 5547.5498046875, 5788.4501953125, 5788.4501953125, 1505561, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(421, 45502, 5811.9501953125, 5858.7001953125, 5480, 5509.89990234375, 5509.89990234375, 1876282, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(422, 45509, 5258, 5587, 5257.0498046875, 5373.5498046875, 5373.5498046875, 2211255, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(423, 45516, 5373.5498046875, 5600, 5302.7998046875, 5563.75, 5563.75, 1115109, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(424, 45523, 5567.9501953125, 5769, 5567.9501953125, 5641.60009765625, 5641.60009765625, 1330168, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # A: Datetime value, formatted the same way as the rest of column A
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value  = $row[2]   # Open
    $ws.Cells.Item($r, 3).Value  = $row[3]   # High
    $ws.Cells.Item($r, 4).Value  = $row[4]   # Low
    $ws.Cells.Item($r, 5).Value  = $row[5]   # Close
    $ws.Cells.Item($r, 6).Value  = $row[6]   # Adj Close
    $ws.Cells.Item($r, 7).Value  = $row[7]   # Volume
    $ws.Cells.Item($r, 8).Value  = $row[8]   # Year
    $ws.Cells.Item($r, 9).Value  = $row[9]   # Month
    $ws.Cells.Item($r, 10).Value = $row[10]  # Day
    $ws.Cells.Item($r, 11).Value = $row[11]  # Hour
    $ws.Cells.Item($r, 12).Value = $row[12]  # Minute
    $ws.Cells.Item($r, 13).Value = $row[13]  # Second
    $ws.Cells.Item($r, 14).Value = $row[14]  # Week
    $ws.Cells.Item($r, 15).Value = $row[15]  # isPivot
    $ws.Cells.Item($r, 16).Value = $row[16]  # two_line_structure
    $ws.Cells.Item($r, 17).Value = $row[17]  # detect_structure
    # Column R (backup) is intentionally left blank for the new rows,
    # matching the not-yet-backed-up state R415/R416 had before this edit.
}
